$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'double[,]' 4,23
$data[0,0] = 0.00220426157237326
$data[0,1] = 0.0191036002939015
$data[0,2] = 0.0337986774430566
$data[0,3] = 0.0352681851579721
$data[0,4] = 0.974283614988979
$data[0,5] = 0.00440852314474651
$data[0,6] = 0.880969875091844
$data[0,7] = 0.00587803085966201
$data[0,8] = 0.141072740631888
$data[0,9] = 0.0705363703159442
$data[0,10] = 0.0675973548861132
$data[0,11] = 0.0183688464364438
$data[0,12] = 0.97722263041881
$data[0,13] = 0.11756061719324
$data[0,14] = 0.0014695077149155
$data[0,15] = 0.000734753857457752
$data[0,16] = 0.157237325495959
$data[0,17] = 0.00440852314474651
$data[0,18] = 0.863335782512858
$data[0,19] = 0.0661278471711976
$data[0,20] = 0.0014695077149155
$data[0,21] = 0.806759735488611
$data[0,22] = 0.0014695077149155
$data[1,0] = 0.979426891991183
$data[1,1] = 0.861131520940485
$data[1,2] = 0.847171197648788
$data[1,3] = 0.017634092578986
$data[1,4] = 0.0161645848640705
$data[1,5] = 0
$data[1,6] = 0.000734753857457752
$data[1,7] = 0.000734753857457752
$data[1,8] = 0.849375459221161
$data[1,9] = 0.000734753857457752
$data[1,10] = 0.000734753857457752
$data[1,11] = 0.938280675973549
$data[1,12] = 0.0014695077149155
$data[1,13] = 0.0360029390154298
$data[1,14] = 0.879500367376929
$data[1,15] = 0.863335782512858
$data[1,16] = 0.000734753857457752
$data[1,17] = 0.00734753857457752
$data[1,18] = 0.0014695077149155
$data[1,19] = 0.0183688464364438
$data[1,20] = 0.000734753857457752
$data[1,21] = 0.0168993387215283
$data[1,22] = 0.83982365907421
$data[2,0] = 0.017634092578986
$data[2,1] = 0
$data[2,2] = 0.000734753857457752
$data[2,3] = 0.946362968405584
$data[2,4] = 0.00955180014695077
$data[2,5] = 0.995591476855254
$data[2,6] = 0.0014695077149155
$data[2,7] = 0.954445260837619
$data[2,8] = 0.00220426157237326
$data[2,9] = 0.000734753857457752
$data[2,10] = 0.815576781778104
$data[2,11] = 0.0014695077149155
$data[2,12] = 0.017634092578986
$data[2,13] = 0.00293901542983101
$data[2,14] = 0.118295371050698
$data[2,15] = 0.0014695077149155
$data[2,16] = 0.831006612784717
$data[2,17] = 0.98677443056576
$data[2,18] = 0.0183688464364438
$data[2,19] = 0.90521675238795
$data[2,20] = 0.017634092578986
$data[2,21] = 0.118295371050698
$data[2,22] = 0
$data[3,0] = 0.000734753857457752
$data[3,1] = 0.119764878765614
$data[3,2] = 0.118295371050698
$data[3,3] = 0.000734753857457752
$data[3,4] = 0
$data[3,5] = 0
$data[3,6] = 0.116825863335783
$data[3,7] = 0.0389419544452608
$data[3,8] = 0.00734753857457752
$data[3,9] = 0.92799412196914
$data[3,10] = 0.116091109478325
$data[3,11] = 0.0418809698750918
$data[3,12] = 0.00367376928728876
$data[3,13] = 0.843497428361499
$data[3,14] = 0.000734753857457752
$data[3,15] = 0.134459955914769
$data[3,16] = 0.0110213078618663
$data[3,17] = 0.0014695077149155
$data[3,18] = 0.116825863335783
$data[3,19] = 0.0102865540044085
$data[3,20] = 0.980161645848641
$data[3,21] = 0.0580455547391624
$data[3,22] = 0.158706833210874

$ws.Range("B2:X5").Value2 = $data
